$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2022" column (M) is being added to the table, mirroring the
# existing style used by column L (2021) in each affected row.
$ws.Range("L3").Copy($ws.Range("M3"))
$ws.Range("L4").Copy($ws.Range("M4"))
$ws.Range("L6").Copy($ws.Range("M6"))
$ws.Range("L7").Copy($ws.Range("M7"))
$ws.Range("L8").Copy($ws.Range("M8"))

# Fill in the new column's data.
$ws.Range("M4").Value = 2022
$ws.Range("M6").Value = 18
$ws.Range("M7").Value = 6.2
$ws.Range("M8").Value = "-"

# Match the author's recorded selection after the edit.
$ws.Range("N4").Select() | Out-Null
